# Updates the "startup" sheet so the CasesTab query (B2) drops the trailing
# `Cohort` column/line, while the SamplesTab (B3) and FilesTab (B4) query
# text stay as-is. Also restores the scroll/selection to B2 and fixes the
# row-2 wrap height now that the text is two lines shorter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- B2: CasesTab query -> drop the trailing Cohort line -------------------
$casesQuery = $ws.Range("B2").Value2

$cohortSuffix = ",`n        coalesce(co.cohort_description, '') AS ``Cohort``"
if ($casesQuery.EndsWith($cohortSuffix)) {
    $casesQuery = $casesQuery.Substring(0, $casesQuery.Length - $cohortSuffix.Length)
}

$ws.Range("B2").Value2 = $casesQuery

# --- Row 2 is two lines shorter now; match the new wrap height -------------
$ws.Rows(2).RowHeight = 304.5

# --- Restore view/selection to B2 (was C4:E4 before the edit) --------------
$ws.Range("B2").Select()
